$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Fn1"
$ws.Cells.Item(2, 3).Value2 = "Tshr"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 19.95578266666667
$ws.Cells.Item(2, 8).Value2 = 59.867348
$ws.Cells.Item(2, 9).Value2 = 0.0117373419656925
$ws.Cells.Item(2, 10).Value2 = 0.0117373419656925
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 0.22859
$ws.Cells.Item(2, 14).Value2 = 0.68577
$ws.Cells.Item(2, 15).Value2 = 0.06730352972305123
$ws.Cells.Item(2, 16).Value2 = 0.06730352972305123
$ws.Cells.Item(2, 17).Value2 = 4.561692359773334
$ws.Cells.Item(2, 18).Value2 = 41.05523123796
$ws.Cells.Item(2, 19).Value2 = 0.0007899645438576015
$ws.Cells.Item(2, 20).Value2 = 0.0007899645438576015

# Row 3
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Fn1"
$ws.Cells.Item(3, 3).Value2 = "Tshr"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 19.95578266666667
$ws.Cells.Item(3, 8).Value2 = 59.867348
$ws.Cells.Item(3, 9).Value2 = 0.0117373419656925
$ws.Cells.Item(3, 10).Value2 = 0.0117373419656925
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 2.180983666666667
$ws.Cells.Item(3, 14).Value2 = 6.542951
$ws.Cells.Item(3, 15).Value2 = 0.6421448840062525
$ws.Cells.Item(3, 16).Value2 = 0.6421448840062525
$ws.Cells.Item(3, 17).Value2 = 43.52323605154978
$ws.Cells.Item(3, 18).Value2 = 391.709124463948
$ws.Cells.Item(3, 19).Value2 = 0.007537074095101327
$ws.Cells.Item(3, 20).Value2 = 0.007537074095101327

# Row 4
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Fn1"
$ws.Cells.Item(4, 3).Value2 = "Tshr"
$ws.Cells.Item(4, 4).Value2 = "M2"
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 19.95578266666667
$ws.Cells.Item(4, 8).Value2 = 59.867348
$ws.Cells.Item(4, 9).Value2 = 0.0117373419656925
$ws.Cells.Item(4, 10).Value2 = 0.0117373419656925
$ws.Cells.Item(4, 11).Value2 = 2
$ws.Cells.Item(4, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(4, 13).Value2 = 0.173517
$ws.Cells.Item(4, 14).Value2 = 0.520551
$ws.Cells.Item(4, 15).Value2 = 0.05108844029465279
$ws.Cells.Item(4, 16).Value2 = 0.05108844029465278
$ws.Cells.Item(4, 17).Value2 = 3.462667540972
$ws.Cells.Item(4, 18).Value2 = 31.164007868748
$ws.Cells.Item(4, 19).Value2 = 0.0005996424942322037
$ws.Cells.Item(4, 20).Value2 = 0.0005996424942322036

# Row 5
$ws.Cells.Item(5, 1).Value2 = "ECs"
$ws.Cells.Item(5, 2).Value2 = "Fn1"
$ws.Cells.Item(5, 3).Value2 = "Tshr"
$ws.Cells.Item(5, 4).Value2 = "sCs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 19.95578266666667
$ws.Cells.Item(5, 8).Value2 = 59.867348
$ws.Cells.Item(5, 9).Value2 = 0.0117373419656925
$ws.Cells.Item(5, 10).Value2 = 0.0117373419656925
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 0.8133136666666667
$ws.Cells.Item(5, 14).Value2 = 2.439941
$ws.Cells.Item(5, 15).Value2 = 0.2394631459760435
$ws.Cells.Item(5, 16).Value2 = 0.2394631459760435
$ws.Cells.Item(5, 17).Value2 = 16.23031077182978
$ws.Cells.Item(5, 18).Value2 = 146.072796946468
$ws.Cells.Item(5, 19).Value2 = 0.002810660832501364
$ws.Cells.Item(5, 20).Value2 = 0.002810660832501363

# Row 6
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Fn1"
$ws.Cells.Item(6, 3).Value2 = "Tshr"
$ws.Cells.Item(6, 4).Value2 = "ECs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 1637.343343333333
$ws.Cells.Item(6, 8).Value2 = 4912.03003
$ws.Cells.Item(6, 9).Value2 = 0.9630320723052701
$ws.Cells.Item(6, 10).Value2 = 0.9630320723052702
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 0.22859
$ws.Cells.Item(6, 14).Value2 = 0.68577
$ws.Cells.Item(6, 15).Value2 = 0.06730352972305123
$ws.Cells.Item(6, 16).Value2 = 0.06730352972305123
$ws.Cells.Item(6, 17).Value2 = 374.2803148525666
$ws.Cells.Item(6, 18).Value2 = 3368.5228336731
$ws.Cells.Item(6, 19).Value2 = 0.06481545770264936
$ws.Cells.Item(6, 20).Value2 = 0.06481545770264938

# Row 7
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Fn1"
$ws.Cells.Item(7, 3).Value2 = "Tshr"
$ws.Cells.Item(7, 4).Value2 = "FAPs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 1637.343343333333
$ws.Cells.Item(7, 8).Value2 = 4912.03003
$ws.Cells.Item(7, 9).Value2 = 0.9630320723052701
$ws.Cells.Item(7, 10).Value2 = 0.9630320723052702
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 2.180983666666667
$ws.Cells.Item(7, 14).Value2 = 6.542951
$ws.Cells.Item(7, 15).Value2 = 0.6421448840062525
$ws.Cells.Item(7, 16).Value2 = 0.6421448840062525
$ws.Cells.Item(7, 17).Value2 = 3571.019088535392
$ws.Cells.Item(7, 18).Value2 = 32139.17179681853
$ws.Cells.Item(7, 19).Value2 = 0.6184061183647687
$ws.Cells.Item(7, 20).Value2 = 0.6184061183647687

# Row 8
$ws.Cells.Item(8, 1).Value2 = "FAPs"
$ws.Cells.Item(8, 2).Value2 = "Fn1"
$ws.Cells.Item(8, 3).Value2 = "Tshr"
$ws.Cells.Item(8, 4).Value2 = "M2"
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 1637.343343333333
$ws.Cells.Item(8, 8).Value2 = 4912.03003
$ws.Cells.Item(8, 9).Value2 = 0.9630320723052701
$ws.Cells.Item(8, 10).Value2 = 0.9630320723052702
$ws.Cells.Item(8, 11).Value2 = 2
$ws.Cells.Item(8, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(8, 13).Value2 = 0.173517
$ws.Cells.Item(8, 14).Value2 = 0.520551
$ws.Cells.Item(8, 15).Value2 = 0.05108844029465279
$ws.Cells.Item(8, 16).Value2 = 0.05108844029465278
$ws.Cells.Item(8, 17).Value2 = 284.10690490517
$ws.Cells.Item(8, 18).Value2 = 2556.96214414653
$ws.Cells.Item(8, 19).Value2 = 0.04919980652780354
$ws.Cells.Item(8, 20).Value2 = 0.04919980652780354

# Row 9
$ws.Cells.Item(9, 1).Value2 = "FAPs"
$ws.Cells.Item(9, 2).Value2 = "Fn1"
$ws.Cells.Item(9, 3).Value2 = "Tshr"
$ws.Cells.Item(9, 4).Value2 = "sCs"
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 1637.343343333333
$ws.Cells.Item(9, 8).Value2 = 4912.03003
$ws.Cells.Item(9, 9).Value2 = 0.9630320723052701
$ws.Cells.Item(9, 10).Value2 = 0.9630320723052702
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 0.8133136666666667
$ws.Cells.Item(9, 14).Value2 = 2.439941
$ws.Cells.Item(9, 15).Value2 = 0.2394631459760435
$ws.Cells.Item(9, 16).Value2 = 0.2394631459760435
$ws.Cells.Item(9, 17).Value2 = 1331.673718158692
$ws.Cells.Item(9, 18).Value2 = 11985.06346342823
$ws.Cells.Item(9, 19).Value2 = 0.2306106897100486
$ws.Cells.Item(9, 20).Value2 = 0.2306106897100486

# Row 10
$ws.Cells.Item(10, 1).Value2 = "M2"
$ws.Cells.Item(10, 2).Value2 = "Fn1"
$ws.Cells.Item(10, 3).Value2 = "Tshr"
$ws.Cells.Item(10, 4).Value2 = "ECs"
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 17.50081933333334
$ws.Cells.Item(10, 8).Value2 = 52.502458
$ws.Cells.Item(10, 9).Value2 = 0.01029341242216722
$ws.Cells.Item(10, 10).Value2 = 0.01029341242216722
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 0.22859
$ws.Cells.Item(10, 14).Value2 = 0.68577
$ws.Cells.Item(10, 15).Value2 = 0.06730352972305123
$ws.Cells.Item(10, 16).Value2 = 0.06730352972305123
$ws.Cells.Item(10, 17).Value2 = 4.000512291406667
$ws.Cells.Item(10, 18).Value2 = 36.00461062266
$ws.Cells.Item(10, 19).Value2 = 0.0006927829889069561
$ws.Cells.Item(10, 20).Value2 = 0.0006927829889069562

# Row 11
$ws.Cells.Item(11, 1).Value2 = "M2"
$ws.Cells.Item(11, 2).Value2 = "Fn1"
$ws.Cells.Item(11, 3).Value2 = "Tshr"
$ws.Cells.Item(11, 4).Value2 = "FAPs"
$ws.Cells.Item(11, 5).Value2 = 3
$ws.Cells.Item(11, 6).Value2 = 1
$ws.Cells.Item(11, 7).Value2 = 17.50081933333334
$ws.Cells.Item(11, 8).Value2 = 52.502458
$ws.Cells.Item(11, 9).Value2 = 0.01029341242216722
$ws.Cells.Item(11, 10).Value2 = 0.01029341242216722
$ws.Cells.Item(11, 11).Value2 = 3
$ws.Cells.Item(11, 12).Value2 = 1
$ws.Cells.Item(11, 13).Value2 = 2.180983666666667
$ws.Cells.Item(11, 14).Value2 = 6.542951
$ws.Cells.Item(11, 15).Value2 = 0.6421448840062525
$ws.Cells.Item(11, 16).Value2 = 0.6421448840062525
$ws.Cells.Item(11, 17).Value2 = 38.16900111928422
$ws.Cells.Item(11, 18).Value2 = 343.521010073558
$ws.Cells.Item(11, 19).Value2 = 0.006609862125861087
$ws.Cells.Item(11, 20).Value2 = 0.006609862125861087

# Row 12
$ws.Cells.Item(12, 1).Value2 = "M2"
$ws.Cells.Item(12, 2).Value2 = "Fn1"
$ws.Cells.Item(12, 3).Value2 = "Tshr"
$ws.Cells.Item(12, 4).Value2 = "M2"
$ws.Cells.Item(12, 5).Value2 = 3
$ws.Cells.Item(12, 6).Value2 = 1
$ws.Cells.Item(12, 7).Value2 = 17.50081933333334
$ws.Cells.Item(12, 8).Value2 = 52.502458
$ws.Cells.Item(12, 9).Value2 = 0.01029341242216722
$ws.Cells.Item(12, 10).Value2 = 0.01029341242216722
$ws.Cells.Item(12, 11).Value2 = 2
$ws.Cells.Item(12, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(12, 13).Value2 = 0.173517
$ws.Cells.Item(12, 14).Value2 = 0.520551
$ws.Cells.Item(12, 15).Value2 = 0.05108844029465279
$ws.Cells.Item(12, 16).Value2 = 0.05108844029465278
$ws.Cells.Item(12, 17).Value2 = 3.036689668262
$ws.Cells.Item(12, 18).Value2 = 27.330207014358
$ws.Cells.Item(12, 19).Value2 = 0.0005258743859581273
$ws.Cells.Item(12, 20).Value2 = 0.0005258743859581273

# Row 13
$ws.Cells.Item(13, 1).Value2 = "M2"
$ws.Cells.Item(13, 2).Value2 = "Fn1"
$ws.Cells.Item(13, 3).Value2 = "Tshr"
$ws.Cells.Item(13, 4).Value2 = "sCs"
$ws.Cells.Item(13, 5).Value2 = 3
$ws.Cells.Item(13, 6).Value2 = 1
$ws.Cells.Item(13, 7).Value2 = 17.50081933333334
$ws.Cells.Item(13, 8).Value2 = 52.502458
$ws.Cells.Item(13, 9).Value2 = 0.01029341242216722
$ws.Cells.Item(13, 10).Value2 = 0.01029341242216722
$ws.Cells.Item(13, 11).Value2 = 3
$ws.Cells.Item(13, 12).Value2 = 1
$ws.Cells.Item(13, 13).Value2 = 0.8133136666666667
$ws.Cells.Item(13, 14).Value2 = 2.439941
$ws.Cells.Item(13, 15).Value2 = 0.2394631459760435
$ws.Cells.Item(13, 16).Value2 = 0.2394631459760435
$ws.Cells.Item(13, 17).Value2 = 14.23365554166423
$ws.Cells.Item(13, 18).Value2 = 128.102899874978
$ws.Cells.Item(13, 19).Value2 = 0.002464892921441048
$ws.Cells.Item(13, 20).Value2 = 0.002464892921441048

# Row 14
$ws.Cells.Item(14, 1).Value2 = "sCs"
$ws.Cells.Item(14, 2).Value2 = "Fn1"
$ws.Cells.Item(14, 3).Value2 = "Tshr"
$ws.Cells.Item(14, 4).Value2 = "ECs"
$ws.Cells.Item(14, 5).Value2 = 3
$ws.Cells.Item(14, 6).Value2 = 1
$ws.Cells.Item(14, 7).Value2 = 25.39612333333333
$ws.Cells.Item(14, 8).Value2 = 76.18836999999999
$ws.Cells.Item(14, 9).Value2 = 0.01493717330687017
$ws.Cells.Item(14, 10).Value2 = 0.01493717330687017
$ws.Cells.Item(14, 11).Value2 = 3
$ws.Cells.Item(14, 12).Value2 = 1
$ws.Cells.Item(14, 13).Value2 = 0.22859
$ws.Cells.Item(14, 14).Value2 = 0.68577
$ws.Cells.Item(14, 15).Value2 = 0.06730352972305123
$ws.Cells.Item(14, 16).Value2 = 0.06730352972305123
$ws.Cells.Item(14, 17).Value2 = 5.805299832766666
$ws.Cells.Item(14, 18).Value2 = 52.24769849489999
$ws.Cells.Item(14, 19).Value2 = 0.001005324487637304
$ws.Cells.Item(14, 20).Value2 = 0.001005324487637304

# Row 15
$ws.Cells.Item(15, 1).Value2 = "sCs"
$ws.Cells.Item(15, 2).Value2 = "Fn1"
$ws.Cells.Item(15, 3).Value2 = "Tshr"
$ws.Cells.Item(15, 4).Value2 = "FAPs"
$ws.Cells.Item(15, 5).Value2 = 3
$ws.Cells.Item(15, 6).Value2 = 1
$ws.Cells.Item(15, 7).Value2 = 25.39612333333333
$ws.Cells.Item(15, 8).Value2 = 76.18836999999999
$ws.Cells.Item(15, 9).Value2 = 0.01493717330687017
$ws.Cells.Item(15, 10).Value2 = 0.01493717330687017
$ws.Cells.Item(15, 11).Value2 = 3
$ws.Cells.Item(15, 12).Value2 = 1
$ws.Cells.Item(15, 13).Value2 = 2.180983666666667
$ws.Cells.Item(15, 14).Value2 = 6.542951
$ws.Cells.Item(15, 15).Value2 = 0.6421448840062525
$ws.Cells.Item(15, 16).Value2 = 0.6421448840062525
$ws.Cells.Item(15, 17).Value2 = 55.38853018665222
$ws.Cells.Item(15, 18).Value2 = 498.49677167987
$ws.Cells.Item(15, 19).Value2 = 0.009591829420521435
$ws.Cells.Item(15, 20).Value2 = 0.009591829420521435

# Row 16
$ws.Cells.Item(16, 1).Value2 = "sCs"
$ws.Cells.Item(16, 2).Value2 = "Fn1"
$ws.Cells.Item(16, 3).Value2 = "Tshr"
$ws.Cells.Item(16, 4).Value2 = "M2"
$ws.Cells.Item(16, 5).Value2 = 3
$ws.Cells.Item(16, 6).Value2 = 1
$ws.Cells.Item(16, 7).Value2 = 25.39612333333333
$ws.Cells.Item(16, 8).Value2 = 76.18836999999999
$ws.Cells.Item(16, 9).Value2 = 0.01493717330687017
$ws.Cells.Item(16, 10).Value2 = 0.01493717330687017
$ws.Cells.Item(16, 11).Value2 = 2
$ws.Cells.Item(16, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(16, 13).Value2 = 0.173517
$ws.Cells.Item(16, 14).Value2 = 0.520551
$ws.Cells.Item(16, 15).Value2 = 0.05108844029465279
$ws.Cells.Item(16, 16).Value2 = 0.05108844029465278
$ws.Cells.Item(16, 17).Value2 = 4.40665913243
$ws.Cells.Item(16, 18).Value2 = 39.65993219186999
$ws.Cells.Item(16, 19).Value2 = 0.0007631168866589179
$ws.Cells.Item(16, 20).Value2 = 0.0007631168866589178

# Row 17
$ws.Cells.Item(17, 1).Value2 = "sCs"
$ws.Cells.Item(17, 2).Value2 = "Fn1"
$ws.Cells.Item(17, 3).Value2 = "Tshr"
$ws.Cells.Item(17, 4).Value2 = "sCs"
$ws.Cells.Item(17, 5).Value2 = 3
$ws.Cells.Item(17, 6).Value2 = 1
$ws.Cells.Item(17, 7).Value2 = 25.39612333333333
$ws.Cells.Item(17, 8).Value2 = 76.18836999999999
$ws.Cells.Item(17, 9).Value2 = 0.01493717330687017
$ws.Cells.Item(17, 10).Value2 = 0.01493717330687017
$ws.Cells.Item(17, 11).Value2 = 3
$ws.Cells.Item(17, 12).Value2 = 1
$ws.Cells.Item(17, 13).Value2 = 0.8133136666666667
$ws.Cells.Item(17, 14).Value2 = 2.439941
$ws.Cells.Item(17, 15).Value2 = 0.2394631459760435
$ws.Cells.Item(17, 16).Value2 = 0.2394631459760435
$ws.Cells.Item(17, 17).Value2 = 20.65501418735222
$ws.Cells.Item(17, 18).Value2 = 185.89512768617
$ws.Cells.Item(17, 19).Value2 = 0.003576902512052512
$ws.Cells.Item(17, 20).Value2 = 0.003576902512052511

